# calibration bug fixes and neucode td
#
# 1. Correct the RLA1_YEAST sequence text (typo fix: S-T-ES... -> S-C-ES...)
#    in the three rows that reference it (rows 2-4, column E).
# 2. Correct the "Observed Precursor Mass" values in column Q for the same
#    rows (calibration fix).
# 3. Widen column E (Sequence) and set column J to the default width so the
#    long sequence text and the newly-visible Modification Codes column are
#    readable.
# 4. Update the saved sheet view so it is scrolled to column F with the
#    active selection on Q2 (reflecting where the analyst was working).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSequence = "STESALSYAALILADSEIEISSEKLLTLTNAANVPVENIWADIFAKALDGQNLKDLLVNFSAGAAAPAGVAGGVAGGEAGEAEAEKEEEEAKEESDDDMGFGLFD"
$newSequence = "SCESALSYAALILADSEIEISSEKLLTLTNAANVPVENIWADIFAKALDGQNLKDLLVNFSAGAAAPAGVAGGVAGGEAGEAEAEKEEEEAKEESDDDMGFGLFD"

$ws.Range("E2").Value = $newSequence
$ws.Range("E3").Value = $newSequence
$ws.Range("E4").Value = $newSequence

$ws.Range("Q2").Value = 10894.13
$ws.Range("Q3").Value = 10894.13
$ws.Range("Q4").Value = 10894.13

$ws.Columns.Item(5).ColumnWidth = 194.7109375
$ws.Columns.Item(10).ColumnWidth = 9.140625

$ws.Range("Q2").Select
$excel.ActiveWindow.ScrollColumn = 6
